# Generate Report for Handback
# Adds a new handed-back file (58c406d2-52f5-4ec7-b0d8-eb99af746003) as
# row 4 of the "Overview", "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$fileId       = "58c406d2-52f5-4ec7-b0d8-eb99af746003"
$mdName       = "$fileId.md"
$zhHash       = "e7684d12687cb2fa5d6398e3af85113c8d02a7cf"
$xlfZh        = "$fileId.$zhHash.zh-cn.xlf"
$xlfDe        = "$fileId.$zhHash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$dateFmt      = "yyyy-mm-dd HH:mm:ss"

$hlColor = 15570276   # RGB(0x64,0x95,0xED) -> matches the workbook's HyperLink style

function Set-HyperlinkStyle($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hlColor
}

# ---------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync
Set-HyperlinkStyle $wsOverview.Range("A4")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/$mdName",
    "",
    "",
    $mdName
)

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Range("C4").Value = $xlfZh
$wsZh.Range("D4").Value = "2016-03-08 02:08:05"
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $xlfZh
$wsZh.Range("G4").Value = "2016-03-08 02:08:43"
$wsZh.Range("H4").Value = "Include"

$wsZh.Range("D4").NumberFormat = $dateFmt
$wsZh.Range("G4").NumberFormat = $dateFmt

Set-HyperlinkStyle $wsZh.Range("A4")
Set-HyperlinkStyle $wsZh.Range("C4")
Set-HyperlinkStyle $wsZh.Range("E4")
Set-HyperlinkStyle $wsZh.Range("F4")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh",
    "",
    "",
    $xlfZh
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh",
    "",
    "",
    $xlfZh
)

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Range("C4").Value = $xlfDe
$wsDe.Range("D4").Value = "2016-03-08 02:08:13"
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $xlfDe
$wsDe.Range("G4").Value = "2016-03-08 02:08:55"
$wsDe.Range("H4").Value = "Include"

$wsDe.Range("D4").NumberFormat = $dateFmt
$wsDe.Range("G4").NumberFormat = $dateFmt

Set-HyperlinkStyle $wsDe.Range("A4")
Set-HyperlinkStyle $wsDe.Range("C4")
Set-HyperlinkStyle $wsDe.Range("E4")
Set-HyperlinkStyle $wsDe.Range("F4")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe",
    "",
    "",
    $xlfDe
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe",
    "",
    "",
    $xlfDe
)

Write-Host "Handback report row added for $fileId"
